$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.435.74"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.14%  '

$ws.Range('D3').Value = "'1.899.48"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.61%  '

$ws.Range('D4').Value = "'1.006"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.48%  '

$ws.Range('D5').Value = "'325.72"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.19%  '

$ws.Range('E6').Value = '  +0.30%  '

$ws.Range('D7').Value = "'0.4784"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.37%  '

$ws.Range('D8').Value = "'0.4054"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.15%  '

$ws.Range('D9').Value = "'0.08066"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.32%  '

$ws.Range('D10').Value = "'1.001"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.17%  '

$ws.Range('D11').Value = "'23.34"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.44%  '

$ws.Range('D12').Value = "'1.937.24"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.75%  '

$ws.Range('D13').Value = "'5.953"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.34%  '

$ws.Range('D14').Value = "'7.070"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.61%  '

$ws.Range('D15').Value = "'90.16"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.31%  '

$ws.Range('D16').Value = "'1.006"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.36%  '

$ws.Range('D17').Value = "'0.06718"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.02%  '

$ws.Range('D18').Value = "'0.00001031"
$ws.Range('D18').Style = 'Normal'

$ws.Range('D19').Value = "'17.60"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.90%  '

$ws.Range('D20').Value = "'1.004"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.18%  '

$ws.Range('D21').Value = "'29.456.90"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.11%  '

$ws.Range('D22').Value = "'5.540"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.51%  '

$ws.Range('E23').Value = '  +2.40%  '

$ws.Range('D24').Value = "'2.161"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.44%  '

$ws.Range('D25').Value = "'2.158.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.13%  '

$ws.Range('D26').Value = "'154.01"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.50%  '

$ws.Range('D27').Value = "'19.87"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.14%  '

$ws.Range('D28').Value = "'6.086"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.53%  '

$ws.Range('D29').Value = "'2.089"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.55%  '

$ws.Range('D30').Value = "'118.35"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.75%  '

$ws.Range('D31').Value = "'1.030"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.35%  '

$ws.Range('D32').Value = "'0.09477"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.13%  '

$ws.Range('D33').Value = "'5.480"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.36%  '

$ws.Range('D34').Value = "'3.548"
$ws.Range('D34').Style = 'Normal'

$ws.Range('D35').Value = "'1.387"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.75%  '

$ws.Range('D36').Value = "'0.06075"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.67%  '

$ws.Range('D37').Value = "'0.02249"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.71%  '

$ws.Range('D38').Value = "'1.170"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.65%  '

$ws.Range('D39').Value = "'0.5878"
$ws.Range('D39').Style = 'Normal'

$ws.Range('D40').Value = "'7.920"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.91%  '

$ws.Range('D41').Value = "'0.1841"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.07%  '

$ws.Range('D42').Value = "'10.22"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.19%  '

$ws.Range('D43').Value = "'1.289"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.75%  '

$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = "'0.07797"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.86%  '

$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = "'2.389"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.65%  '

$ws.Range('D46').Value = "'12.24"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.27%  '

$ws.Range('D47').Value = "'0.5525"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.80%  '

$ws.Range('D48').Value = "'1.920"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.44%  '

$ws.Range('D49').Value = "'114.02"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.66%  '

$ws.Range('D50').Value = "'72.35"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.21%  '

$ws.Range('D51').Value = "'0.2927"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.16%  '
